$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    1  = "Troian Landscape -005"
    2  = "----------------"
    3  = "01: 088(058): StingRat x3, Treant   x2"
    4  = "02: 089(059): Panther  x1, StingRat x3"
    5  = "03: 090(05A): Cannibal x1, Treant   x2"
    6  = "04: 091(05B): Python   x1, StingRat x2"
    7  = "05: 088(058): StingRat x3, Treant   x2"
    8  = "06: 089(059): Panther  x1, StingRat x3"
    9  = "07: 090(05A): Cannibal x1, Treant   x2"
    10 = "08: 091(05B): Python   x1, StingRat x2"
    12 = "Agart Island -006"
    13 = "------------"
    14 = "01: 136(088): Roc      x1, Roc Baby x2"
    15 = "02: 139(08B): HugeCell x3"
    16 = "03: 210(0D2): FlameDog x2, BlackLiz x2"
    17 = "04: 145(091): Ironback x2, BlackLiz x2"
    18 = "05: 137(089): Roc      x1, Roc Baby x3"
    19 = "06: 138(08A): HugeCell x4"
    20 = "07: 210(0D2): FlameDog x2, BlackLiz x2"
    21 = "08: 145(091): Ironback x2, BlackLiz x2"
    23 = "Eblan Cave Entrance -042"
    24 = "-------------------"
    25 = "01: 140(08C): GiantBat x3"
    26 = "02: 143(08F): Ironback x2"
    27 = "03: 146(092): Skull    x3"
    28 = "04: 145(091): Ironback x2, BlackLiz x2"
    29 = "05: 147(093): Skull    x4"
    30 = "06: 142(08E): GiantBat x3, Cave Bat x3"
    31 = "07: 150(096): Staleman x1, Skull    x2"
    32 = "08: 149(095): Staleman x1"
}

foreach ($row in $values.Keys | Sort-Object) {
    $ws.Cells.Item($row, 10).Value = $values[$row]
}

$ws.Range("J23:J32").Select()
